$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "per question" marking scheme (row 11)
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Update the totals row (row 12), driven by the new marking scheme
$ws.Range("B12").Value = 90
$ws.Range("C12").Value = -7.199999999999999
$ws.Range("E12").Value = "82.8/140"
